$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2
$ws.Range("G3").Formula = "=(E3-F3)"

$ws.Range("F4").Value = 1
$ws.Range("G4").Formula = "=(E4-F4)"

$ws.Range("F5").Value = 3
$ws.Range("G5").Formula = "=(E5-F5)"

$ws.Range("D23").Select()
